$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2934.6897
$ws.Range("I98").Value = 1400.5652
$ws.Range("J98").Value = 8815.5
$ws.Range("K98").Value = 1400.5652
$ws.Range("L98").Value = 8815.5
$ws.Range("M98").Value = 97.4348
$ws.Range("N98").Value = -11811.5

$ws.Range("H122").Value = 2934.6897
$ws.Range("I122").Value = 1400.5652
$ws.Range("J122").Value = 8815.5
$ws.Range("K122").Value = 4201.6956
$ws.Range("L122").Value = 26446.5
$ws.Range("M122").Value = -1751.6956
$ws.Range("N122").Value = -31346.5

$ws.Range("H132").Value = 91024.22
$ws.Range("I132").Value = 106607.734
$ws.Range("J132").Value = 6428
$ws.Range("K132").Value = 319823.202
$ws.Range("L132").Value = 19284
$ws.Range("M132").Value = -317293.202
$ws.Range("N132").Value = -24344

$ws.Range("H141").Value = 11699.091
$ws.Range("I141").Value = 12369
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 37107
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -31927
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1885.2812
$ws.Range("I61").Value = 1306
$ws.Range("J61").Value = 3954.1428
$ws.Range("K61").Value = 1306
$ws.Range("L61").Value = 3954.1428
$ws.Range("M61").Value = -1094
$ws.Range("N61").Value = -4378.1428

$ws.Range("H74").Value = 1095.1333
$ws.Range("I74").Value = 809.587
$ws.Range("J74").Value = 2033.3572
$ws.Range("K74").Value = 809.587
$ws.Range("L74").Value = 2033.3572
$ws.Range("M74").Value = 64.41300000000001
$ws.Range("N74").Value = -3781.3572

$ws.Range("H77").Value = 1095.1333
$ws.Range("I77").Value = 809.587
$ws.Range("J77").Value = 2033.3572
$ws.Range("K77").Value = 4047.935
$ws.Range("L77").Value = 10166.786
$ws.Range("M77").Value = 320.0650000000001
$ws.Range("N77").Value = -18902.786

$ws.Range("H122").Value = 2551.4092
$ws.Range("I122").Value = 1667
$ws.Range("K122").Value = 5001
$ws.Range("M122").Value = -2551

$ws.Range("H124").Value = 29429
$ws.Range("J124").Value = 29429
$ws.Range("L124").Value = 29429
$ws.Range("N124").Value = -39249

$ws.Range("H125").Value = 42739.375
$ws.Range("J125").Value = 42739.375
$ws.Range("L125").Value = 42739.375
$ws.Range("N125").Value = -52579.375

$ws.Range("H136").Value = 1885.2812
$ws.Range("I136").Value = 1306
$ws.Range("J136").Value = 3954.1428
$ws.Range("K136").Value = 3918
$ws.Range("L136").Value = 11862.4284
$ws.Range("M136").Value = -1368
$ws.Range("N136").Value = -16962.4284

$ws.Range("H137").Value = 41780
$ws.Range("J137").Value = 41780
$ws.Range("L137").Value = 41780
$ws.Range("N137").Value = -51980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H55").Value = 19000
$ws.Range("J55").Value = 19000
$ws.Range("L55").Value = 19000
$ws.Range("N55").Value = -19546

$ws.Range("H59").Value = 118858
$ws.Range("J59").Value = 118858
$ws.Range("L59").Value = 118858
$ws.Range("N59").Value = -120552

$ws.Range("H87").Value = 40800
$ws.Range("J87").Value = 40800
$ws.Range("L87").Value = 40800
$ws.Range("N87").Value = -43296

$ws.Range("H90").Value = 40800
$ws.Range("J90").Value = 40800
$ws.Range("L90").Value = 122400
$ws.Range("N90").Value = -134880

$ws.Range("H137").Value = 35513.332
$ws.Range("J137").Value = 40770
$ws.Range("L137").Value = 40770
$ws.Range("N137").Value = -50970

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 5405
$ws.Range("I32").Value = 540
$ws.Range("K32").Value = 540
$ws.Range("M32").Value = -224

$ws.Range("H45").Value = 11533.5
$ws.Range("I45").Value = 3067
$ws.Range("K45").Value = 3067
$ws.Range("M45").Value = -2474

$ws.Range("H58").Value = 2226.842
$ws.Range("I58").Value = 1083.125
$ws.Range("J58").Value = 8326.666999999999
$ws.Range("K58").Value = 1083.125
$ws.Range("L58").Value = 8326.666999999999
$ws.Range("M58").Value = -880.125
$ws.Range("N58").Value = -8732.666999999999

$ws.Range("H132").Value = 3550.25
$ws.Range("I132").Value = 2918.652
$ws.Range("J132").Value = 6455.6
$ws.Range("K132").Value = 8755.956
$ws.Range("L132").Value = 19366.8
$ws.Range("M132").Value = -6225.956
$ws.Range("N132").Value = -24426.8

$ws.Range("H136").Value = 2226.842
$ws.Range("I136").Value = 1083.125
$ws.Range("J136").Value = 8326.666999999999
$ws.Range("K136").Value = 3249.375
$ws.Range("L136").Value = 24980.001
$ws.Range("M136").Value = -699.375
$ws.Range("N136").Value = -30080.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 89.5
$ws.Range("J12").Value = 108.833336
$ws.Range("L12").Value = 326.500008
$ws.Range("N12").Value = -672.500008

$ws.Range("H106").Value = 4196.6665
$ws.Range("J106").Value = 4196.6665
$ws.Range("L106").Value = 12589.9995
$ws.Range("N106").Value = -14481.9995

$ws.Range("H113").Value = 565.8125
$ws.Range("I113").Value = 543.5714
$ws.Range("J113").Value = 596.95
$ws.Range("K113").Value = 1630.7142
$ws.Range("L113").Value = 1790.85
$ws.Range("M113").Value = 539.2857999999999
$ws.Range("N113").Value = -6130.85

$ws.Range("H131").Value = 8929557
$ws.Range("J131").Value = 995.67926
$ws.Range("L131").Value = 2987.03778
$ws.Range("N131").Value = -13067.03778

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2908.4814
$ws.Range("I122").Value = 1682
$ws.Range("J122").Value = 4692.4546
$ws.Range("K122").Value = 5046
$ws.Range("L122").Value = 14077.3638
$ws.Range("M122").Value = -2596
$ws.Range("N122").Value = -18977.3638

$ws.Range("H132").Value = 3211.5789
$ws.Range("I132").Value = 1414.4
$ws.Range("K132").Value = 4243.200000000001
$ws.Range("M132").Value = -1713.200000000001

$ws.Range("H137").Value = 48350
$ws.Range("J137").Value = 48350
$ws.Range("L137").Value = 48350
$ws.Range("N137").Value = -58550

$ws.Range("H140").Value = 39716
$ws.Range("J140").Value = 39716
$ws.Range("L140").Value = 39716
$ws.Range("N140").Value = -50076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2369.75
$ws.Range("I16").Value = 2369.75
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2369.75
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2199.75
$ws.Range("N16").ClearContents()

$ws.Range("H32").Value = 5571
$ws.Range("I32").Value = 856.5
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 856.5
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -539.5
$ws.Range("N32").Value = -15634

$ws.Range("H127").Value = 63000
$ws.Range("J127").Value = 63000
$ws.Range("L127").Value = 63000
$ws.Range("N127").Value = -72920

$ws.Range("H133").Value = 34750
$ws.Range("J133").Value = 34750
$ws.Range("L133").Value = 34750
$ws.Range("N133").Value = -39810

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4560.706
$ws.Range("I122").Value = 2553.2
$ws.Range("K122").Value = 7659.599999999999
$ws.Range("M122").Value = -5209.599999999999

$ws.Range("H132").Value = 10421121
$ws.Range("I132").Value = 5320.2383
$ws.Range("J132").Value = 30305832
$ws.Range("K132").Value = 15960.7149
$ws.Range("L132").Value = 90917496
$ws.Range("M132").Value = -13430.7149
$ws.Range("N132").Value = -90922556

$ws.Range("H136").Value = 1119.6945
$ws.Range("I136").Value = 510.68967
$ws.Range("J136").Value = 3642.7144
$ws.Range("K136").Value = 1532.06901
$ws.Range("L136").Value = 10928.1432
$ws.Range("M136").Value = 1017.93099
$ws.Range("N136").Value = -16028.1432

